$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest snapshot.
# For the Price column, force a Text number format before writing so Excel keeps
# the exact original string (e.g. "82.00", "1.00", "61.251.20") instead of
# reinterpreting / reformatting it as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.251.20"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.355.84"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "403.29"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.75"
$ws.Range("E6").Value = "  +14.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.600"
$ws.Range("E7").Value = "  +5.76%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.668"
$ws.Range("E9").Value = "  +7.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  +9.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.69"
$ws.Range("E11").Value = "  +6.88%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.930.15"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.45"
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.54"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.342.86"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.251.32"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  +4.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.02"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000131"
$ws.Range("E20").Value = "  +15.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.23"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "82.00"
$ws.Range("E22").Value = "  +11.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.01"
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "304.77"
$ws.Range("E24").Value = "  +3.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.48"
$ws.Range("E26").Value = "  +13.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.67"
$ws.Range("E27").Value = "  +9.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.40"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.43"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.173"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.115"
$ws.Range("E31").Value = "  +2.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.65"
$ws.Range("E32").Value = "  +4.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.59"
$ws.Range("E33").Value = "  +4.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.81"
$ws.Range("E34").Value = "  +5.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0482"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.30"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.99"
$ws.Range("E41").Value = "  +6.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.124"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.57"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.93"
$ws.Range("E44").Value = "  +4.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.88"
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.282"
$ws.Range("E46").Value = "  -3.46%  "
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.51"
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.136.38"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.703.61"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  -0.53%  "
